$wb = $excel.ActiveWorkbook

# The new "Italy" test-data sheet was created by duplicating the existing
# "Germany" sheet (same layout/market label) and then updating the part
# number in B4. Before duplicating, all cells on "Germany" were selected
# (Ctrl+A), which is why that sheet's saved selection ends up as the full
# sheet range.
$germany = $wb.Worksheets.Item("Germany")
$germany.Cells.Select() | Out-Null

# Copy "Germany" to the end of the workbook (after the last existing sheet,
# "Slovakia") to create the new sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$germany.Copy($null, $lastSheet)

$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Update the part-number cell for the new market.
$italy.Range("B4").Value = "NGC-3475/T1820/NGC-3145/T2446"

# Leave the selection on B4, on the new (now active) "Italy" tab.
$italy.Range("B4").Select() | Out-Null
